$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "Snow can now transition in and out" -- the "Controls" mark (row 17) moves
# from 7.5/10 up to full marks (10/10), and is now annotated DONE (matching
# the style already used on the other completed rows) in column E instead of
# carrying the leftover "speed" scratch note in column F.
$ws.Range("D17").Value = 10
$ws.Range("F17").ClearContents()
$ws.Range("E17").Value = "DONE"
$ws.Range("E17").Font.Color = $ws.Range("E16").Font.Color

# Move the active selection in the frozen bottom-right pane to D21.
$ws.Range("D21").Select()
